# "Add BOL in stage" - refresh the randomized QA fixture SKU values in the
# ManageProducts "Input" sheet (column B, rows 2-7) with newly generated
# product codes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value2 = "prodRuXj"
$ws.Range("B3").Value2 = "prodVCFD"
$ws.Range("B4").Value2 = "prodUHkc"
$ws.Range("B5").Value2 = "prodbwCv"
$ws.Range("B6").Value2 = "prodtLdr"
$ws.Range("B7").Value2 = "prodyVIP"
